# InscribirFacturas.xlsx - "refactorización features y data"
#
# The "Datos" sheet holds a single data-driven test row (row 2). This
# refresh swaps the test user/credentials/reference values used by the
# automated scenario for a new set (new robot user + matching
# numeroDocumento), keeping the rest of the row intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# numeroDocumento
$ws.Range("B2").Value = 48646663

# usuario / clave (new autotest credentials replace the old robot ones)
$ws.Range("D2").Value = "autotest11"
$ws.Range("E2").Value = "1234"

# Land the selection back on the edited cell (also scrolls the view back
# to column A, dropping the previous C1 scroll position)
$ws.Range("B2").Select()
